$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.91741
$ws.Range("H2").Value = 20.75223
$ws.Range("I2").Value = 0.2334435312127427
$ws.Range("J2").Value = 0.2334435312127427
$ws.Range("M2").Value = 1.016190333333333
$ws.Range("N2").Value = 3.048571
$ws.Range("O2").Value = 0.1011512660469171
$ws.Range("P2").Value = 0.1011512660469171
$ws.Range("Q2").Value = 7.029405173703333
$ws.Range("R2").Value = 63.26464656333
$ws.Range("S2").Value = 0.02361310873263193
$ws.Range("T2").Value = 0.02361310873263193
$ws.Range("G3").Value = 6.91741
$ws.Range("H3").Value = 20.75223
$ws.Range("I3").Value = 0.2334435312127427
$ws.Range("J3").Value = 0.2334435312127427
$ws.Range("O3").Value = 0.1112440260843082
$ws.Range("P3").Value = 0.1112440260843082
$ws.Range("Q3").Value = 7.730791349046667
$ws.Range("R3").Value = 69.57712214142001
$ws.Range("S3").Value = 0.02596919827544336
$ws.Range("T3").Value = 0.02596919827544336
$ws.Range("G4").Value = 6.91741
$ws.Range("H4").Value = 20.75223
$ws.Range("I4").Value = 0.2334435312127427
$ws.Range("J4").Value = 0.2334435312127427
$ws.Range("M4").Value = 5.77911
$ws.Range("N4").Value = 17.33733
$ws.Range("O4").Value = 0.5752507910667645
$ws.Range("P4").Value = 0.5752507910667645
$ws.Range("Q4").Value = 39.9764733051
$ws.Range("R4").Value = 359.7882597459001
$ws.Range("S4").Value = 0.1342885759995492
$ws.Range("T4").Value = 0.1342885759995492
$ws.Range("G5").Value = 6.91741
$ws.Range("H5").Value = 20.75223
$ws.Range("I5").Value = 0.2334435312127427
$ws.Range("J5").Value = 0.2334435312127427
$ws.Range("M5").Value = 0.3634723333333333
$ws.Range("N5").Value = 1.090417
$ws.Range("O5").Value = 0.03617992169743831
$ws.Range("P5").Value = 0.03617992169743831
$ws.Range("Q5").Value = 2.514287153323334
$ws.Range("R5").Value = 22.62858437991
$ws.Range("S5").Value = 0.008445968680050526
$ws.Range("T5").Value = 0.008445968680050526
$ws.Range("G6").Value = 6.91741
$ws.Range("H6").Value = 20.75223
$ws.Range("I6").Value = 0.2334435312127427
$ws.Range("J6").Value = 0.2334435312127427
$ws.Range("M6").Value = 1.769887
$ws.Range("N6").Value = 5.309661
$ws.Range("O6").Value = 0.1761739951045719
$ws.Range("P6").Value = 0.1761739951045719
$ws.Range("Q6").Value = 12.24303403267
$ws.Range("R6").Value = 110.18730629403
$ws.Range("S6").Value = 0.04112667952506771
$ws.Range("T6").Value = 0.04112667952506771
$ws.Range("I7").Value = 0.2633623201546029
$ws.Range("J7").Value = 0.2633623201546028
$ws.Range("M7").Value = 1.016190333333333
$ws.Range("N7").Value = 3.048571
$ws.Range("O7").Value = 0.1011512660469171
$ws.Range("P7").Value = 0.1011512660469171
$ws.Range("Q7").Value = 7.930313794671665
$ws.Range("R7").Value = 71.37282415204498
$ws.Range("S7").Value = 0.02663943211269159
$ws.Range("T7").Value = 0.02663943211269159
$ws.Range("I8").Value = 0.2633623201546029
$ws.Range("J8").Value = 0.2633623201546028
$ws.Range("O8").Value = 0.1112440260843082
$ws.Range("P8").Value = 0.1112440260843082
$ws.Range("S8").Value = 0.02929748481290256
$ws.Range("T8").Value = 0.02929748481290256
$ws.Range("I9").Value = 0.2633623201546029
$ws.Range("J9").Value = 0.2633623201546028
$ws.Range("M9").Value = 5.77911
$ws.Range("N9").Value = 17.33733
$ws.Range("O9").Value = 0.5752507910667645
$ws.Range("P9").Value = 0.5752507910667645
$ws.Range("Q9").Value = 45.09997217114999
$ws.Range("R9").Value = 405.89974954035
$ws.Range("S9").Value = 0.1514993830061138
$ws.Range("T9").Value = 0.1514993830061138
$ws.Range("I10").Value = 0.2633623201546029
$ws.Range("J10").Value = 0.2633623201546028
$ws.Range("M10").Value = 0.3634723333333333
$ws.Range("N10").Value = 1.090417
$ws.Range("O10").Value = 0.03617992169743831
$ws.Range("P10").Value = 0.03617992169743831
$ws.Range("Q10").Value = 2.836525367801666
$ws.Range("R10").Value = 25.528728310215
$ws.Range("S10").Value = 0.009528428121249211
$ws.Range("T10").Value = 0.009528428121249207
$ws.Range("I11").Value = 0.2633623201546029
$ws.Range("J11").Value = 0.2633623201546028
$ws.Range("M11").Value = 1.769887
$ws.Range("N11").Value = 5.309661
$ws.Range("O11").Value = 0.1761739951045719
$ws.Range("P11").Value = 0.1761739951045719
$ws.Range("Q11").Value = 13.812136201955
$ws.Range("R11").Value = 124.309225817595
$ws.Range("S11").Value = 0.04639759210164571
$ws.Range("T11").Value = 0.04639759210164569
$ws.Range("G12").Value = 6.430676666666667
$ws.Range("H12").Value = 19.29203
$ws.Range("I12").Value = 0.2170176220802376
$ws.Range("J12").Value = 0.2170176220802376
$ws.Range("M12").Value = 1.016190333333333
$ws.Range("N12").Value = 3.048571
$ws.Range("O12").Value = 0.1011512660469171
$ws.Range("P12").Value = 0.1011512660469171
$ws.Range("Q12").Value = 6.534791465458889
$ws.Range("R12").Value = 58.81312318913
$ws.Range("S12").Value = 0.02195160722790742
$ws.Range("T12").Value = 0.02195160722790742
$ws.Range("G13").Value = 6.430676666666667
$ws.Range("H13").Value = 19.29203
$ws.Range("I13").Value = 0.2170176220802376
$ws.Range("J13").Value = 0.2170176220802376
$ws.Range("O13").Value = 0.1112440260843082
$ws.Range("P13").Value = 0.1112440260843082
$ws.Range("Q13").Value = 7.186825638957778
$ws.Range("R13").Value = 64.68143075062
$ws.Range("S13").Value = 0.02414191401144848
$ws.Range("T13").Value = 0.02414191401144848
$ws.Range("G14").Value = 6.430676666666667
$ws.Range("H14").Value = 19.29203
$ws.Range("I14").Value = 0.2170176220802376
$ws.Range("J14").Value = 0.2170176220802376
$ws.Range("M14").Value = 5.77911
$ws.Range("N14").Value = 17.33733
$ws.Range("O14").Value = 0.5752507910667645
$ws.Range("P14").Value = 0.5752507910667645
$ws.Range("Q14").Value = 37.1635878311
$ws.Range("R14").Value = 334.4722904799
$ws.Range("S14").Value = 0.1248395587770848
$ws.Range("T14").Value = 0.1248395587770848
$ws.Range("G15").Value = 6.430676666666667
$ws.Range("H15").Value = 19.29203
$ws.Range("I15").Value = 0.2170176220802376
$ws.Range("J15").Value = 0.2170176220802376
$ws.Range("M15").Value = 0.3634723333333333
$ws.Range("N15").Value = 1.090417
$ws.Range("O15").Value = 0.03617992169743831
$ws.Range("P15").Value = 0.03617992169743831
$ws.Range("Q15").Value = 2.337373052945556
$ws.Range("R15").Value = 21.03635747651
$ws.Range("S15").Value = 0.007851680573827253
$ws.Range("T15").Value = 0.007851680573827253
$ws.Range("G16").Value = 6.430676666666667
$ws.Range("H16").Value = 19.29203
$ws.Range("I16").Value = 0.2170176220802376
$ws.Range("J16").Value = 0.2170176220802376
$ws.Range("M16").Value = 1.769887
$ws.Range("N16").Value = 5.309661
$ws.Range("O16").Value = 0.1761739951045719
$ws.Range("P16").Value = 0.1761739951045719
$ws.Range("Q16").Value = 11.38157103353667
$ws.Range("R16").Value = 102.43413930183
$ws.Range("S16").Value = 0.03823286148996961
$ws.Range("T16").Value = 0.03823286148996961
$ws.Range("G17").Value = 4.144241333333333
$ws.Range("H17").Value = 12.432724
$ws.Range("I17").Value = 0.1398567283204463
$ws.Range("J17").Value = 0.1398567283204463
$ws.Range("M17").Value = 1.016190333333333
$ws.Range("N17").Value = 3.048571
$ws.Range("O17").Value = 0.1011512660469171
$ws.Range("P17").Value = 0.1011512660469171
$ws.Range("Q17").Value = 4.211337981933778
$ws.Range("R17").Value = 37.902041837404
$ws.Range("S17").Value = 0.01414668513479287
$ws.Range("T17").Value = 0.01414668513479287
$ws.Range("G18").Value = 4.144241333333333
$ws.Range("H18").Value = 12.432724
$ws.Range("I18").Value = 0.1398567283204463
$ws.Range("J18").Value = 0.1398567283204463
$ws.Range("O18").Value = 0.1112440260843082
$ws.Range("P18").Value = 0.1112440260843082
$ws.Range("Q18").Value = 4.631540569099556
$ws.Range("R18").Value = 41.683865121896
$ws.Range("S18").Value = 0.01555822553334573
$ws.Range("T18").Value = 0.01555822553334573
$ws.Range("G19").Value = 4.144241333333333
$ws.Range("H19").Value = 12.432724
$ws.Range("I19").Value = 0.1398567283204463
$ws.Range("J19").Value = 0.1398567283204463
$ws.Range("M19").Value = 5.77911
$ws.Range("N19").Value = 17.33733
$ws.Range("O19").Value = 0.5752507910667645
$ws.Range("P19").Value = 0.5752507910667645
$ws.Range("Q19").Value = 23.95002653188
$ws.Range("R19").Value = 215.55023878692
$ws.Range("S19").Value = 0.08045269360234629
$ws.Range("T19").Value = 0.08045269360234629
$ws.Range("G20").Value = 4.144241333333333
$ws.Range("H20").Value = 12.432724
$ws.Range("I20").Value = 0.1398567283204463
$ws.Range("J20").Value = 0.1398567283204463
$ws.Range("M20").Value = 0.3634723333333333
$ws.Range("N20").Value = 1.090417
$ws.Range("O20").Value = 0.03617992169743831
$ws.Range("P20").Value = 0.03617992169743831
$ws.Range("Q20").Value = 1.506317067323111
$ws.Range("R20").Value = 13.556853605908
$ws.Range("S20").Value = 0.005060005479493649
$ws.Range("T20").Value = 0.005060005479493649
$ws.Range("G21").Value = 4.144241333333333
$ws.Range("H21").Value = 12.432724
$ws.Range("I21").Value = 0.1398567283204463
$ws.Range("J21").Value = 0.1398567283204463
$ws.Range("M21").Value = 1.769887
$ws.Range("N21").Value = 5.309661
$ws.Range("O21").Value = 0.1761739951045719
$ws.Range("P21").Value = 0.1761739951045719
$ws.Range("Q21").Value = 7.334838860729334
$ws.Range("R21").Value = 66.01354974656401
$ws.Range("S21").Value = 0.02463911857046775
$ws.Range("T21").Value = 0.02463911857046775
$ws.Range("G22").Value = 4.335755333333334
$ws.Range("H22").Value = 13.007266
$ws.Range("I22").Value = 0.1463197982319706
$ws.Range("J22").Value = 0.1463197982319706
$ws.Range("M22").Value = 1.016190333333333
$ws.Range("N22").Value = 3.048571
$ws.Range("O22").Value = 0.1011512660469171
$ws.Range("P22").Value = 0.1011512660469171
$ws.Range("Q22").Value = 4.405952657431778
$ws.Range("R22").Value = 39.653573916886
$ws.Range("S22").Value = 0.01480043283889329
$ws.Range("T22").Value = 0.01480043283889329
$ws.Range("G23").Value = 4.335755333333334
$ws.Range("H23").Value = 13.007266
$ws.Range("I23").Value = 0.1463197982319706
$ws.Range("J23").Value = 0.1463197982319706
$ws.Range("O23").Value = 0.1112440260843082
$ws.Range("P23").Value = 0.1112440260843082
$ws.Range("Q23").Value = 4.845573678951556
$ws.Range("R23").Value = 43.610163110564
$ws.Range("S23").Value = 0.01627720345116805
$ws.Range("T23").Value = 0.01627720345116804
$ws.Range("G24").Value = 4.335755333333334
$ws.Range("H24").Value = 13.007266
$ws.Range("I24").Value = 0.1463197982319706
$ws.Range("J24").Value = 0.1463197982319706
$ws.Range("M24").Value = 5.77911
$ws.Range("N24").Value = 17.33733
$ws.Range("O24").Value = 0.5752507910667645
$ws.Range("P24").Value = 0.5752507910667645
$ws.Range("Q24").Value = 25.05680700442
$ws.Range("R24").Value = 225.51126303978
$ws.Range("S24").Value = 0.08417057968167045
$ws.Range("T24").Value = 0.08417057968167044
$ws.Range("G25").Value = 4.335755333333334
$ws.Range("H25").Value = 13.007266
$ws.Range("I25").Value = 0.1463197982319706
$ws.Range("J25").Value = 0.1463197982319706
$ws.Range("M25").Value = 0.3634723333333333
$ws.Range("N25").Value = 1.090417
$ws.Range("O25").Value = 0.03617992169743831
$ws.Range("P25").Value = 0.03617992169743831
$ws.Range("Q25").Value = 1.575927107769111
$ws.Range("R25").Value = 14.183343969922
$ws.Range("S25").Value = 0.005293838842817668
$ws.Range("T25").Value = 0.005293838842817668
$ws.Range("G26").Value = 4.335755333333334
$ws.Range("H26").Value = 13.007266
$ws.Range("I26").Value = 0.1463197982319706
$ws.Range("J26").Value = 0.1463197982319706
$ws.Range("M26").Value = 1.769887
$ws.Range("N26").Value = 5.309661
$ws.Range("O26").Value = 0.1761739951045719
$ws.Range("P26").Value = 0.1761739951045719
$ws.Range("Q26").Value = 7.673796999647335
$ws.Range("R26").Value = 69.06417299682602
$ws.Range("S26").Value = 0.02577774341742114
$ws.Range("T26").Value = 0.02577774341742113

Write-Output "Applied 278 cell updates"
